$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country name pairs (swap adjacent shared-string rows) ---
$ws.Range("A98").Value = "Albania"
$ws.Range("A99").Value = "Finlandia"

$ws.Range("A129").Value = "Sudan del Sur"
$ws.Range("A130").Value = "Eslovenia"

$ws.Range("A144").Value = "Jordania"
$ws.Range("A145").Value = "Malta"

$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Update "last updated" timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 18:35"

# --- Update statistic values (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 5666985
$ws.Range("C4").Value = 11011
$ws.Range("D4").Value = 3013595
$ws.Range("E4").Value = 2477900
$ws.Range("G4").Value = 416
$ws.Range("H4").Value = 175490

$ws.Range("B5").Value = 3418306
$ws.Range("C5").Value = 6434
$ws.Range("E5").Value = 753956
$ws.Range("G5").Value = 152
$ws.Range("H5").Value = 110171

$ws.Range("B6").Value = 2823078
$ws.Range("C6").Value = 56452
$ws.Range("D6").Value = 2083198
$ws.Range("E6").Value = 686042
$ws.Range("G6").Value = 824
$ws.Range("H6").Value = 53838

$ws.Range("B15").Value = 321098
$ws.Range("C15").Value = 812
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 41397

$ws.Range("B43").Value = 69801
$ws.Range("C43").Value = 128
$ws.Range("D43").Value = 67647
$ws.Range("E43").Value = 1532
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 622

$ws.Range("D49").Value = 52810
$ws.Range("E49").Value = 3194

$ws.Range("E53").Value = 3484
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 176

$ws.Range("B66").Value = 31015
$ws.Range("C66").Value = 379
$ws.Range("D66").Value = 17612
$ws.Range("E66").Value = 12897
$ws.Range("G66").Value = 19
$ws.Range("H66").Value = 506

$ws.Range("B74").Value = 20686
$ws.Range("C74").Value = 203
$ws.Range("D74").Value = 15615
$ws.Range("E74").Value = 4667
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 404

$ws.Range("B94").Value = 8711
$ws.Range("C94").Value = 54
$ws.Range("D94").Value = 8112
$ws.Range("E94").Value = 545
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 54

$ws.Range("B98").Value = 7812
$ws.Range("C98").Value = 158
$ws.Range("D98").Value = 3928
$ws.Range("E98").Value = 3650
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 234

$ws.Range("B99").Value = 7805
$ws.Range("C99").Value = 29
$ws.Range("D99").Value = 7100
$ws.Range("E99").Value = 371
$ws.Range("H99").Value = 334

$ws.Range("B129").Value = 2494
$ws.Range("C129").Value = 4
$ws.Range("D129").Value = 1175
$ws.Range("E129").Value = 1272
$ws.Range("H129").Value = 47

$ws.Range("B130").Value = 2493
$ws.Range("C130").Value = 37
$ws.Range("D130").Value = 2079
$ws.Range("E130").Value = 285
$ws.Range("H130").Value = 129

$ws.Range("B136").Value = 2035
$ws.Range("C136").Value = 8
$ws.Range("D136").Value = 1903

$ws.Range("B144").Value = 1482
$ws.Range("C144").Value = 44
$ws.Range("D144").Value = 1259
$ws.Range("E144").Value = 212
$ws.Range("H144").Value = 11

$ws.Range("B145").Value = 1470
$ws.Range("C145").Value = 47
$ws.Range("D145").Value = 784
$ws.Range("E145").Value = 677
$ws.Range("H145").Value = 9

$ws.Range("B157").Value = 1024
$ws.Range("C157").Value = 19
$ws.Range("D157").Value = 875
$ws.Range("E157").Value = 96

$ws.Range("B165").Value = 704
$ws.Range("C165").Value = 2
$ws.Range("E165").Value = 5

$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

